# Saldo_guide.xlsx refresh: new export snapshot (Dt. Referencia 24/07 -> 26/07),
# updated Saldo Previsto / Vl. Total figures for the clients whose balance moved,
# one client renamed to reflect an estate ("ESPOLIO ..."), and the workbook/sheet
# renamed to match the new export timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the client who is now listed as an estate (ESPOLIO) ---
$ws.Range("C145").Value = "ESPOLIO AGNES LOUIZE MOURA DE SANTANA"

# --- Reference date (column G) moved from 2024-07-24 (45497) to 2024-07-26 (45499) for every data row ---
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = 45499
}

# --- Updated Saldo Previsto (E) / Vl. Total (H) amounts for the rows that changed ---
$updates = @{
    5   = 592.38
    8   = 415.56
    15  = 817.98
    17  = 416.81
    19  = 2.27
    36  = 22503.3
    43  = 1619.16
    49  = 196.14
    58  = 200.45
    60  = 337.78
    99  = 360.49
    104 = 997.24
    107 = 54894.98
    108 = 1195.6400000000001
    132 = 217.07
    143 = 1358.97
    158 = 697.9
    172 = 867.15
    173 = 1060.76
    231 = 1417.32
    235 = 458.62
    249 = 252.01
    264 = 1895.5
    265 = 912.81
    270 = 554.72
    271 = 778.96
    273 = 788.86
    274 = 22.34
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $ws.Cells.Item($row, 5).Value = $value
    $ws.Cells.Item($row, 8).Value = $value
}

# --- Sheet / export name bumped to the new extraction timestamp ---
$ws.Name = "IClientBalance-20240726-100632-"
